$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-19, row 13 untouched) were rearranged / shuffled.
# Each target row below receives the full record (columns D, H, I, J, K, L, M, N, O, P)
# that originally belonged to a different source row.

$data = @{
    2  = @{ D=44511; H="Sin especificar"; I="Primera"; J=600;  K=1300; L=1400; M=1350; N="`$/kilo";    O="Provincia de Linares"; P=1350 }
    3  = @{ D=44860; H="Sin especificar"; I="Primera"; J=1100; K=1500; L=1700; M=1609; N="`$/kilo";    O="Provincia de Linares"; P=1609 }
    4  = @{ D=44881; H="Sin especificar"; I="Primera"; J=200;  K=2600; L=2700; M=2650; N="`$/kilo";    O="Provincia de Linares"; P=2650 }
    5  = @{ D=44881; H="Sin especificar"; I="Segunda"; J=100;  K=2400; L=2400; M=2400; N="`$/kilo";    O="Provincia de Linares"; P=2400 }
    6  = @{ D=44489; H="Sin especificar"; I="Primera"; J=600;  K=1400; L=1500; M=1450; N="`$/kilo";    O="Provincia de Linares"; P=1450 }
    7  = @{ D=44477; H="Sin especificar"; I="Primera"; J=500;  K=1400; L=1500; M=1460; N="`$/kilo";    O="Provincia de Linares"; P=1460 }
    8  = @{ D=44839; H="Sin especificar"; I="Primera"; J=500;  K=1700; L=1800; M=1760; N="`$/kilo";    O="Provincia de Linares"; P=1760 }
    9  = @{ D=44876; H="Sin especificar"; I="Primera"; J=350;  K=1500; L=1600; M=1557; N="`$/kilo";    O="Provincia de Linares"; P=1557 }
    10 = @{ D=44510; H="Sin especificar"; I="Primera"; J=600;  K=1300; L=1400; M=1350; N="`$/kilo";    O="Provincia de Linares"; P=1350 }
    11 = @{ D=44875; H="Sin especificar"; I="Primera"; J=300;  K=1500; L=1600; M=1550; N="`$/kilo";    O="Provincia de Linares"; P=1550 }
    12 = @{ D=44468; H="Verde";           I="Primera"; J=500;  K=1800; L=2000; M=1920; N="`$/kilo";    O="Provincia de Linares"; P=1920 }
    14 = @{ D=44496; H="Sin especificar"; I="Primera"; J=550;  K=1500; L=2000; M=1773; N="`$/paquete"; O="Provincia de Linares"; P=1773 }
    15 = @{ D=44524; H="Sin especificar"; I="Primera"; J=200;  K=1500; L=1600; M=1550; N="`$/kilo";    O="Provincia de Talca";   P=1550 }
    16 = @{ D=44868; H="Sin especificar"; I="Primera"; J=1000; K=1200; L=1300; M=1250; N="`$/kilo";    O="Región del Maule";     P=1250 }
    17 = @{ D=44868; H="Sin especificar"; I="Segunda"; J=200;  K=1000; L=1000; M=1000; N="`$/kilo";    O="Región del Maule";     P=1000 }
    18 = @{ D=44519; H="Sin especificar"; I="Primera"; J=250;  K=1200; L=1300; M=1240; N="`$/kilo";    O="Provincia de Linares"; P=1240 }
    19 = @{ D=44526; H="Sin especificar"; I="Primera"; J=100;  K=1500; L=1600; M=1550; N="`$/kilo";    O="Provincia de Linares"; P=1550 }
}

foreach ($row in $data.Keys) {
    $rec = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $rec.D   # D: Fecha
    $ws.Cells.Item($row, 8).Value  = $rec.H   # H: Variedad
    $ws.Cells.Item($row, 9).Value  = $rec.I   # I: Calidad
    $ws.Cells.Item($row, 10).Value = $rec.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $rec.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $rec.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $rec.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $rec.N   # N: Unidad de comercializacion
    $ws.Cells.Item($row, 15).Value = $rec.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $rec.P   # P: Precio $/Kg
}
